# Dry Run Fix 3
# Update spiciness price values on the "Menu" sheet and set the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu")

# Update spice level prices (column H)
$ws.Range("H10").Value = 0.1
$ws.Range("H12").Value = 0.3
$ws.Range("H13").Value = 0.4
$ws.Range("H14").Value = 0.5

# Set the active cell/selection as recorded in the saved view state
$ws.Activate()
$ws.Range("H16").Select()
